# Product backlog update 04/30/2018
# Renumber the "Story Priority" column (D) on rows 7-24 to account for a
# newly-inserted backlog item, shifting H/M/L priority numbers up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value  = "H, 1"
$ws.Range("D8").Value  = "H, 2"
$ws.Range("D9").Value  = "H, 3"
$ws.Range("D10").Value = "H, 4"
$ws.Range("D11").Value = "H, 5"
$ws.Range("D12").Value = "M, 6"
$ws.Range("D13").Value = "M, 7"
$ws.Range("D14").Value = "M, 8"
$ws.Range("D15").Value = "M, 9"
$ws.Range("D16").Value = "M, 10"
$ws.Range("D17").Value = "L, 11"
$ws.Range("D18").Value = "L, 12"
$ws.Range("D19").Value = "L, 13"
$ws.Range("D20").Value = "L, 14"
$ws.Range("D21").Value = "L, 15"
$ws.Range("D22").Value = "L, 16"
$ws.Range("D23").Value = "L, 17"
$ws.Range("D24").Value = "L, 18"

# The author's selection ended on D28 when they saved the file.
$ws.Range("D28").Select()
